$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.560778711498519
$ws.Range("H2").Value = 0.982142857142857
$ws.Range("J2").Value = 5.7
$ws.Range("K2").Value = 0.0249829001367989
$ws.Range("L2").Value = -0.183187990322654
$ws.Range("M2").Value = 0.319729996363022
$ws.Range("N2").Value = 0.438296493628051
$ws.Range("P2").Value = "As likely as not increasing"

$ws.Range("F3").Value = 0.226131977792684
$ws.Range("H3").Value = 0.824561403508772
$ws.Range("J3").Value = 0.045
$ws.Range("K3").Value = 0.0029535040431266
$ws.Range("L3").Value = -0.0019800812101678
$ws.Range("M3").Value = 0.006714372097811
$ws.Range("N3").Value = 6.5633423180593
$ws.Range("P3").Value = "Unlikely improving"

$ws.Range("F4").Value = 0.983311243253705
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = -91.230655657192
$ws.Range("L4").Value = -200.655681068827
$ws.Range("M4").Value = -21.8272317755256
$ws.Range("N4").Value = -16.5873919376713
$ws.Range("P4").Value = "Extremely likely improving"

$ws.Range("F5").Value = 0.236351531311268
$ws.Range("J5").Value = 0.0622215874544633
$ws.Range("K5").Value = 0.0047420822083413
$ws.Range("L5").Value = -0.0050612958497988
$ws.Range("M5").Value = 0.0152178128846579
$ws.Range("N5").Value = 7.62128129857149
$ws.Range("P5").Value = "Unlikely improving"

$ws.Range("D6").Value = $false
$ws.Range("F6").Value = 0.192515854989737
$ws.Range("H6").Value = 0.614035087719298
$ws.Range("K6").Value = 0.001003434065934
$ws.Range("L6").Value = -0.0007692403438804
$ws.Range("M6").Value = 0.0035664844151136
$ws.Range("N6").Value = 5.28123192596877
$ws.Range("P6").Value = "Unlikely improving"

$ws.Range("F7").Value = 0.184880240211656
$ws.Range("H7").Value = 0.964912280701754
$ws.Range("J7").Value = 0.251
$ws.Range("K7").Value = 0.0160549450549451
$ws.Range("L7").Value = -0.0254478386784755
$ws.Range("M7").Value = 0.0326946213107602
$ws.Range("N7").Value = 6.39639245216935
$ws.Range("P7").Value = "Unlikely improving"

$ws.Range("F8").Value = 0.002809441074419
$ws.Range("H8").Value = 0.696428571428571
$ws.Range("J8").Value = 7.21
$ws.Range("K8").Value = -0.0492250673854446
$ws.Range("L8").Value = -0.0790081032101419
$ws.Range("M8").Value = -0.0272948778384987
$ws.Range("N8").Value = -0.682733250838344
$ws.Range("P8").Value = "Exceptionally unlikely increasing"

$ws.Range("F9").Value = 0.673089756838225
$ws.Range("H9").Value = 0.947368421052632
$ws.Range("J9").Value = 0.487
$ws.Range("K9").Value = -0.0240824175824176
$ws.Range("L9").Value = -0.101900242238031
$ws.Range("M9").Value = 0.0375340602886444
$ws.Range("N9").Value = -4.94505494505494
$ws.Range("P9").Value = "Likely improving"

$ws.Range("B10").Value = "Total Nitrogen"
$ws.Range("F10").Value = 0.552043600863612
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.842105263157895
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1.26
$ws.Range("K10").Value = -0.0061407434402332
$ws.Range("L10").Value = -0.12341875807846
$ws.Range("M10").Value = 0.115905106683097
$ws.Range("N10").Value = -0.487360590494702
$ws.Range("P10").Value = "As likely as not improving"
$ws.Range("W10").Value = "g/m3"

$ws.Range("B11").Value = "Total Phosphorus"
$ws.Range("D11").Value = $true
$ws.Range("F11").Value = 0.382479927714759
$ws.Range("H11").Value = 0.912280701754386
$ws.Range("J11").Value = 0.139
$ws.Range("K11").Value = 0.0055188873626373
$ws.Range("L11").Value = -0.0084910863695685
$ws.Range("M11").Value = 0.0140480769230769
$ws.Range("N11").Value = 3.97042256304846
$ws.Range("P11").Value = "As likely as not improving"

$ws.Range("B12").Value = "Turbidity"
$ws.Range("F12").Value = 0.5
$ws.Range("H12").Value = 0.982456140350877
$ws.Range("J12").Value = 4.1
$ws.Range("K12").Value = 0.0200686813186813
$ws.Range("L12").Value = -0.51188292608263
$ws.Range("M12").Value = 0.277799157924156
$ws.Range("N12").Value = 0.489480032162959
$ws.Range("P12").Value = "As likely as not improving"
$ws.Range("W12").Value = "NTU/FNU"

$ws.Range("B13").Value = "Visual Clarity"
$ws.Range("C13").Value = 10
$ws.Range("F13").Value = 0.427723999081393
$ws.Range("H13").Value = 0.795454545454545
$ws.Range("J13").Value = 0.335
$ws.Range("K13").Value = -0.002813550420168
$ws.Range("L13").Value = -0.0541009309566701
$ws.Range("M13").Value = 0.0242937274504725
$ws.Range("N13").Value = -0.839865797065095
$ws.Range("P13").Value = "As likely as not improving"
$ws.Range("W13").Value = "m"

$ws.Range("B14").Value = "Dissolved Oxygen Concentration"
$ws.Range("D14").Value = $true
$ws.Range("F14").Value = 0.854743991218061
$ws.Range("H14").Value = 0.945945945945946
$ws.Range("J14").Value = 5.97
$ws.Range("K14").Value = 0.0535164835164834
$ws.Range("L14").Value = -0.0288753451939663
$ws.Range("M14").Value = 0.160000790139708
$ws.Range("N14").Value = 0.896423509488834
$ws.Range("P14").Value = "Likely increasing"
$ws.Range("W14").Value = "g/m3"

$ws.Range("B15").Value = "Dissolved Reactive Phosphorus"
$ws.Range("F15").Value = 0.999995350094737
$ws.Range("H15").Value = 0.786324786324786
$ws.Range("J15").Value = 0.059
$ws.Range("K15").Value = -0.0073187727019949
$ws.Range("L15").Value = -0.0138399744855887
$ws.Range("M15").Value = -0.0035951594161541
$ws.Range("N15").Value = -12.4046994949067
$ws.Range("P15").Value = "Virtually certain improving"
$ws.Range("W15").Value = "mg/L"

$ws.Range("B16").Value = "E. coli"
$ws.Range("D16").Value = $false
$ws.Range("F16").Value = 0.999976827289679
$ws.Range("H16").Value = 0.777777777777778
$ws.Range("J16").Value = 870
$ws.Range("K16").Value = -102.195965833232
$ws.Range("L16").Value = -159.138862433627
$ws.Range("M16").Value = -61.4512511253774
$ws.Range("N16").Value = -11.746662739452
$ws.Range("W16").Value = "E. coli/100 mL"

$ws.Range("B17").Value = "Ammoniacal Nitrogen (NH4)"
$ws.Range("F17").Value = 0.576832660031419
$ws.Range("G17").Value = 0.055045871559633
$ws.Range("H17").Value = 0.954128440366973
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 0.0747490526712699
$ws.Range("K17").Value = -0.000469304979891
$ws.Range("L17").Value = -0.0047438932077628
$ws.Range("M17").Value = 0.0037896928014854
$ws.Range("N17").Value = -0.627840705828074
$ws.Range("P17").Value = "As likely as not improving"
$ws.Range("W17").Value = "mg/L"

$ws.Range("B18").Value = "Nitrite Nitrogen (NO2)"
$ws.Range("F18").Value = 0.752140606763398
$ws.Range("G18").Value = 0.0170940170940171
$ws.Range("H18").Value = 0.47008547008547
$ws.Range("J18").Value = 0.02
$ws.Range("K18").Value = -0.0003529558359552
$ws.Range("L18").Value = -0.0012722895933915
$ws.Range("M18").Value = 0.0004626744106291
$ws.Range("N18").Value = -1.76477917977636
$ws.Range("P18").Value = "Likely improving"

$ws.Range("B19").Value = "Nitrate Nitrogen (NO3)"
$ws.Range("F19").Value = 0.278719637612447
$ws.Range("G19").Value = 0.0256410256410256
$ws.Range("H19").Value = 0.94017094017094
$ws.Range("J19").Value = 0.2925
$ws.Range("K19").Value = 0.0040009128251939
$ws.Range("L19").Value = -0.0089296285708767
$ws.Range("M19").Value = 0.0139019289461326
$ws.Range("N19").Value = 1.36783344451076
$ws.Range("P19").Value = "Unlikely improving"

$ws.Range("B20").Value = "pH"
$ws.Range("F20").Value = 0.012353285661632
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0.552631578947368
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 7.23
$ws.Range("K20").Value = -0.0100343406593404
$ws.Range("L20").Value = -0.0199414249702019
$ws.Range("M20").Value = -0.0050133144823866
$ws.Range("N20").Value = -0.138787560986728
$ws.Range("P20").Value = "Extremely unlikely increasing"
$ws.Range("W20").Value = ""

$ws.Range("B21").Value = "SIN (Soluble Inorganic nitrogen)"
$ws.Range("F21").Value = 0.879661952119208
$ws.Range("H21").Value = 0.94017094017094
$ws.Range("J21").Value = 1.07
$ws.Range("K21").Value = -0.0240248041432252
$ws.Range("L21").Value = -0.0690955347383401
$ws.Range("M21").Value = 0.0090634113882962
$ws.Range("N21").Value = -2.24530879843226
$ws.Range("P21").Value = "Likely improving"
$ws.Range("W21").Value = "g/m3"

$ws.Range("B22").Value = "Total Nitrogen"
$ws.Range("D22").Value = $false
$ws.Range("F22").Value = 0.999767990763701
$ws.Range("H22").Value = 0.871794871794872
$ws.Range("J22").Value = 2.1
$ws.Range("K22").Value = -0.149721624558673
$ws.Range("L22").Value = -0.220146316648275
$ws.Range("M22").Value = -0.0750727439897241
$ws.Range("N22").Value = -7.12960116946061
$ws.Range("P22").Value = "Virtually certain improving"

$ws.Range("B23").Value = "Total Phosphorus"
$ws.Range("D23").Value = $true
$ws.Range("E23").Value = "ok"
$ws.Range("F23").Value = 0.999993036006252
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0.914529914529915
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0.208
$ws.Range("K23").Value = -0.0207544815795715
$ws.Range("L23").Value = -0.0327794117769439
$ws.Range("M23").Value = -0.0117567594869892
$ws.Range("N23").Value = -9.97811614402476
$ws.Range("W23").Value = "g/m3"

$ws.Range("B24").Value = "Turbidity"
$ws.Range("F24").Value = 0.999999100971773
$ws.Range("H24").Value = 0.94017094017094
$ws.Range("J24").Value = 5.64
$ws.Range("K24").Value = -0.559183135554541
$ws.Range("L24").Value = -0.789206807818945
$ws.Range("M24").Value = -0.317436113969417
$ws.Range("N24").Value = -9.91459460203088
$ws.Range("W24").Value = "NTU/FNU"

$ws.Range("B25").Value = "Visual Clarity"
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = $false
$ws.Range("F25").Value = 0.001555186415166
$ws.Range("H25").Value = 0.710526315789474
$ws.Range("J25").Value = 0.42
$ws.Range("K25").Value = -0.0738582116252615
$ws.Range("L25").Value = -0.122943382265489
$ws.Range("M25").Value = -0.02740544187078
$ws.Range("N25").Value = -17.5852884822051
$ws.Range("P25").Value = "Exceptionally unlikely improving"
$ws.Range("W25").Value = "m"

$ws.Range("B26").Value = "Dissolved Oxygen Concentration"
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = $true
$ws.Range("F26").Value = 0.00000032193941401637
$ws.Range("H26").Value = 0.923076923076923
$ws.Range("J26").Value = 7.59
$ws.Range("K26").Value = -0.192616641337386
$ws.Range("L26").Value = -0.260892857142857
$ws.Range("M26").Value = -0.133459581264641
$ws.Range("N26").Value = -2.53776866057162
$ws.Range("P26").Value = "Exceptionally unlikely increasing"
$ws.Range("W26").Value = "g/m3"

$ws.Range("B27").Value = "Dissolved Reactive Phosphorus"
$ws.Range("F27").Value = 0.615946477900092
$ws.Range("H27").Value = 0.676136363636364
$ws.Range("J27").Value = 0.0525
$ws.Range("K27").Value = -0.0002504571428571
$ws.Range("L27").Value = -0.0013464443558658
$ws.Range("M27").Value = 0.0012024701542938
$ws.Range("N27").Value = -0.477061224489796
$ws.Range("P27").Value = "As likely as not improving"
$ws.Range("W27").Value = "mg/L"

$ws.Range("B28").Value = "E. coli"
$ws.Range("D28").Value = $false
$ws.Range("F28").Value = 0.992526242556823
$ws.Range("H28").Value = 0.784090909090909
$ws.Range("J28").Value = 855
$ws.Range("K28").Value = -34.4671053973997
$ws.Range("L28").Value = -65.4531415776593
$ws.Range("M28").Value = -10.8042929268797
$ws.Range("N28").Value = -4.0312403973567
$ws.Range("P28").Value = "Virtually certain improving"
$ws.Range("W28").Value = "E. coli/100 mL"

$ws.Range("B29").Value = "Ammoniacal Nitrogen (NH4)"
$ws.Range("F29").Value = 0.0001757458060435
$ws.Range("G29").Value = 0.0958083832335329
$ws.Range("H29").Value = 0.904191616766467
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 0.0622649705531497
$ws.Range("K29").Value = 0.0035667789109907
$ws.Range("L29").Value = 0.0018910596698125
$ws.Range("M29").Value = 0.0055523689365568
$ws.Range("N29").Value = 5.72838769424314
$ws.Range("P29").Value = "Exceptionally unlikely improving"
$ws.Range("W29").Value = "mg/L"

$ws.Range("B30").Value = "Nitrite Nitrogen (NO2)"
$ws.Range("D30").Value = $true
$ws.Range("F30").Value = 0.220169245069128
$ws.Range("G30").Value = 0.0481927710843374
$ws.Range("H30").Value = 0.379518072289157
$ws.Range("I30").Value = 2
$ws.Range("J30").Value = 0.019
$ws.Range("K30").Value = 0.0002496582365003
$ws.Range("L30").Value = -0.0002276822603127
$ws.Range("M30").Value = 0.000751508989199
$ws.Range("N30").Value = 1.31399071842285
$ws.Range("P30").Value = "Unlikely improving"

$ws.Range("B31").Value = "Nitrate Nitrogen (NO3)"
$ws.Range("F31").Value = 0.461661088395773
$ws.Range("G31").Value = 0.0542168674698795
$ws.Range("H31").Value = 0.909638554216867
$ws.Range("J31").Value = 0.334
$ws.Range("K31").Value = 0.000207854199372
$ws.Range("L31").Value = -0.0077320601885345
$ws.Range("M31").Value = 0.0054931493650231
$ws.Range("N31").Value = 0.0622317962191783
$ws.Range("P31").Value = "As likely as not improving"

$ws.Range("B32").Value = "pH"
$ws.Range("E32").Value = "ok"
$ws.Range("F32").Value = 0.084997290041083
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0.47093023255814
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 7.24
$ws.Range("K32").Value = -0.0065106486620392
$ws.Range("L32").Value = -0.0132761236749108
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = -0.0899260864922552
$ws.Range("P32").Value = "Very unlikely increasing"
$ws.Range("W32").Value = ""

$ws.Range("B33").Value = "SIN (Soluble Inorganic nitrogen)"
$ws.Range("F33").Value = 0.161873021713892
$ws.Range("H33").Value = 0.954545454545455
$ws.Range("J33").Value = 0.8098
$ws.Range("K33").Value = 0.0114212972744559
$ws.Range("L33").Value = -0.0059635783510388
$ws.Range("M33").Value = 0.0283701160985274
$ws.Range("N33").Value = 1.4103849437461
$ws.Range("P33").Value = "Unlikely improving"
$ws.Range("W33").Value = "g/m3"

$ws.Range("B34").Value = "Total Nitrogen"
$ws.Range("F34").Value = 0.627105234628898
$ws.Range("H34").Value = 0.840909090909091
$ws.Range("J34").Value = 2.042
$ws.Range("K34").Value = -0.0066895604395604
$ws.Range("L34").Value = -0.0535736486835946
$ws.Range("M34").Value = 0.0180646690389568
$ws.Range("N34").Value = -0.32759845443489
$ws.Range("P34").Value = "As likely as not improving"

$ws.Range("B35").Value = "Total Phosphorus"
$ws.Range("F35").Value = 0.429849430059029
$ws.Range("H35").Value = 0.840909090909091
$ws.Range("J35").Value = 0.1755
$ws.Range("K35").Value = 0.000501717032967
$ws.Range("L35").Value = -0.003442920827875
$ws.Range("M35").Value = 0.0040076205060861
$ws.Range("N35").Value = 0.285878651263267

$ws.Range("B36").Value = "Turbidity"
$ws.Range("D36").Value = $true
$ws.Range("F36").Value = 0.999997995598715
$ws.Range("H36").Value = 0.943181818181818
$ws.Range("J36").Value = 6.405
$ws.Range("K36").Value = -0.305057965594615
$ws.Range("L36").Value = -0.475996578704558
$ws.Range("M36").Value = -0.193722827430725
$ws.Range("N36").Value = -4.76280976728517
$ws.Range("P36").Value = "Virtually certain improving"
$ws.Range("W36").Value = "NTU/FNU"

$ws.Rows.Item(37).Delete()
